# Natmi following Dr Hou advice
# Rewrites the LR-pair data rows (A2:T19) on Sheet1 with updated values:
#  - "M1"/"M2" sending-cluster rows now use pooled stats across 3 cells (was 1)
#  - a new sending cluster "Neutro" is added (rows 14-19)
#  - receptor-expressing-cell counts / derived expression values are recomputed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 18,20

# Row 2: "M1", "Cd5l", "Cd5" -> "ECs"
$data[0,0] = "M1"
$data[0,1] = "Cd5l"
$data[0,2] = "Cd5"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 4.859690000000001
$data[0,7] = 14.57907
$data[0,8] = 0.4445218291063105
$data[0,9] = 0.4445218291063105
$data[0,10] = 2
$data[0,11] = 1
$data[0,12] = 0.5150980000000001
$data[0,13] = 1.030196
$data[0,14] = 0.07656787870625796
$data[0,15] = 0.05529560532417913
$data[0,16] = 2.50321659962
$data[0,17] = 15.01929959772
$data[0,18] = 0.03403609349329591
$data[0,19] = 0.02458010362024475

# Row 3: "M1", "Cd5l", "Cd5" -> "FAPs"
$data[1,0] = "M1"
$data[1,1] = "Cd5l"
$data[1,2] = "Cd5"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 4.859690000000001
$data[1,7] = 14.57907
$data[1,8] = 0.4445218291063105
$data[1,9] = 0.4445218291063105
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 1.497116333333333
$data[1,13] = 4.491349
$data[1,14] = 0.2225421605593971
$data[1,15] = 0.2410724383293534
$data[1,16] = 7.275521273936667
$data[1,17] = 65.47969146542999
$data[1,18] = 0.09892484826513341
$data[1,19] = 0.1071619612332824

# Row 4: "M1", "Cd5l", "Cd5" -> "M1"
$data[2,0] = "M1"
$data[2,1] = "Cd5l"
$data[2,2] = "Cd5"
$data[2,3] = "M1"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 4.859690000000001
$data[2,7] = 14.57907
$data[2,8] = 0.4445218291063105
$data[2,9] = 0.4445218291063105
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.9008219999999999
$data[2,13] = 2.702466
$data[2,14] = 0.1339046737357332
$data[2,15] = 0.1450544297764824
$data[2,16] = 4.37771566518
$data[2,17] = 39.39944098662
$data[2,18] = 0.05952355049489186
$data[2,19] = 0.06447986044421482

# Row 5: "M1", "Cd5l", "Cd5" -> "M2"
$data[3,0] = "M1"
$data[3,1] = "Cd5l"
$data[3,2] = "Cd5"
$data[3,3] = "M2"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 4.859690000000001
$data[3,7] = 14.57907
$data[3,8] = 0.4445218291063105
$data[3,9] = 0.4445218291063105
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 2.363391666666667
$data[3,13] = 7.090175
$data[3,14] = 0.3513115688057693
$data[3,15] = 0.3805640077027689
$data[3,16] = 11.48535084858334
$data[3,17] = 103.36815763725
$data[3,18] = 0.156165661151748
$data[3,19] = 0.1691690087960629

# Row 6: "M1", "Cd5l", "Cd5" -> "Neutro"
$data[4,0] = "M1"
$data[4,1] = "Cd5l"
$data[4,2] = "Cd5"
$data[4,3] = "Neutro"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 4.859690000000001
$data[4,7] = 14.57907
$data[4,8] = 0.4445218291063105
$data[4,9] = 0.4445218291063105
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.414697
$data[4,13] = 1.244091
$data[4,14] = 0.06164355053960423
$data[4,15] = 0.06677638519598536
$data[4,16] = 2.01529886393
$data[4,17] = 18.13768977537
$data[4,18] = 0.02740190383847216
$data[4,19] = 0.02968356088842697

# Row 7: "M1", "Cd5l", "Cd5" -> "sCs"
$data[5,0] = "M1"
$data[5,1] = "Cd5l"
$data[5,2] = "Cd5"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 4.859690000000001
$data[5,7] = 14.57907
$data[5,8] = 0.4445218291063105
$data[5,9] = 0.4445218291063105
$data[5,10] = 2
$data[5,11] = 1
$data[5,12] = 1.036213
$data[5,13] = 2.072426
$data[5,14] = 0.1540301676532382
$data[5,15] = 0.1112371336712308
$data[5,16] = 5.035673953970001
$data[5,17] = 30.21404372382001
$data[5,18] = 0.06846977186276909
$data[5,19] = 0.04944733412407866

# Row 8: "M2", "Cd5l", "Cd5" -> "ECs"
$data[6,0] = "M2"
$data[6,1] = "Cd5l"
$data[6,2] = "Cd5"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 6.046898666666667
$data[6,7] = 18.140696
$data[6,8] = 0.5531172679177431
$data[6,9] = 0.5531172679177431
$data[6,10] = 2
$data[6,11] = 1
$data[6,12] = 0.5150980000000001
$data[6,13] = 1.030196
$data[6,14] = 0.07656787870625796
$data[6,15] = 0.05529560532417913
$data[6,16] = 3.114745409402667
$data[6,17] = 18.688472456416
$data[6,18] = 0.04235101588026254
$data[6,19] = 0.03058495414476777

# Row 9: "M2", "Cd5l", "Cd5" -> "FAPs"
$data[7,0] = "M2"
$data[7,1] = "Cd5l"
$data[7,2] = "Cd5"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 6.046898666666667
$data[7,7] = 18.140696
$data[7,8] = 0.5531172679177431
$data[7,9] = 0.5531172679177431
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 1.497116333333333
$data[7,13] = 4.491349
$data[7,14] = 0.2225421605593971
$data[7,15] = 0.2410724383293534
$data[7,16] = 9.052910759878223
$data[7,17] = 81.476196838904
$data[7,18] = 0.1230919118451254
$data[7,19] = 0.1333413284590006

# Row 10: "M2", "Cd5l", "Cd5" -> "M1"
$data[8,0] = "M2"
$data[8,1] = "Cd5l"
$data[8,2] = "Cd5"
$data[8,3] = "M1"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 6.046898666666667
$data[8,7] = 18.140696
$data[8,8] = 0.5531172679177431
$data[8,9] = 0.5531172679177431
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.9008219999999999
$data[8,13] = 2.702466
$data[8,14] = 0.1339046737357332
$data[8,15] = 0.1450544297764824
$data[8,16] = 5.447179350704
$data[8,17] = 49.024614156336
$data[8,18] = 0.07406498729812551
$data[8,19] = 0.08023210989733406

# Row 11: "M2", "Cd5l", "Cd5" -> "M2"
$data[9,0] = "M2"
$data[9,1] = "Cd5l"
$data[9,2] = "Cd5"
$data[9,3] = "M2"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 6.046898666666667
$data[9,7] = 18.140696
$data[9,8] = 0.5531172679177431
$data[9,9] = 0.5531172679177431
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 2.363391666666667
$data[9,13] = 7.090175
$data[9,14] = 0.3513115688057693
$data[9,15] = 0.3805640077027689
$data[9,16] = 14.29118991797778
$data[9,17] = 128.6207092618
$data[9,18] = 0.1943164951257433
$data[9,19] = 0.2104965242083825

# Row 12: "M2", "Cd5l", "Cd5" -> "Neutro"
$data[10,0] = "M2"
$data[10,1] = "Cd5l"
$data[10,2] = "Cd5"
$data[10,3] = "Neutro"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 6.046898666666667
$data[10,7] = 18.140696
$data[10,8] = 0.5531172679177431
$data[10,9] = 0.5531172679177431
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 0.414697
$data[10,13] = 1.244091
$data[10,14] = 0.06164355053960423
$data[10,15] = 0.06677638519598536
$data[10,16] = 2.507630736370667
$data[10,17] = 22.568676627336
$data[10,18] = 0.03409611225921521
$data[10,19] = 0.03693517174102625

# Row 13: "M2", "Cd5l", "Cd5" -> "sCs"
$data[11,0] = "M2"
$data[11,1] = "Cd5l"
$data[11,2] = "Cd5"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 6.046898666666667
$data[11,7] = 18.140696
$data[11,8] = 0.5531172679177431
$data[11,9] = 0.5531172679177431
$data[11,10] = 2
$data[11,11] = 1
$data[11,12] = 1.036213
$data[11,13] = 2.072426
$data[11,14] = 0.1540301676532382
$data[11,15] = 0.1112371336712308
$data[11,16] = 6.265875008082667
$data[11,17] = 37.59525004849601
$data[11,18] = 0.08519674550927103
$data[11,19] = 0.06152717946723195

# Row 14: "Neutro", "Cd5l", "Cd5" -> "ECs"
$data[12,0] = "Neutro"
$data[12,1] = "Cd5l"
$data[12,2] = "Cd5"
$data[12,3] = "ECs"
$data[12,4] = 1
$data[12,5] = 0.3333333333333333
$data[12,6] = 0.02581033333333333
$data[12,7] = 0.077431
$data[12,8] = 0.002360902975946389
$data[12,9] = 0.002360902975946389
$data[12,10] = 2
$data[12,11] = 1
$data[12,12] = 0.5150980000000001
$data[12,13] = 1.030196
$data[12,14] = 0.07656787870625796
$data[12,15] = 0.05529560532417913
$data[12,16] = 0.01329485107933333
$data[12,17] = 0.07976910647600001
$data[12,18] = 0.0001807693326995066
$data[12,19] = 0.0001305475591666115

# Row 15: "Neutro", "Cd5l", "Cd5" -> "FAPs"
$data[13,0] = "Neutro"
$data[13,1] = "Cd5l"
$data[13,2] = "Cd5"
$data[13,3] = "FAPs"
$data[13,4] = 1
$data[13,5] = 0.3333333333333333
$data[13,6] = 0.02581033333333333
$data[13,7] = 0.077431
$data[13,8] = 0.002360902975946389
$data[13,9] = 0.002360902975946389
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 1.497116333333333
$data[13,13] = 4.491349
$data[13,14] = 0.2225421605593971
$data[13,15] = 0.2410724383293534
$data[13,16] = 0.03864107160211111
$data[13,17] = 0.347769644419
$data[13,18] = 0.0005254004491382196
$data[13,19] = 0.0005691486370704229

# Row 16: "Neutro", "Cd5l", "Cd5" -> "M1"
$data[14,0] = "Neutro"
$data[14,1] = "Cd5l"
$data[14,2] = "Cd5"
$data[14,3] = "M1"
$data[14,4] = 1
$data[14,5] = 0.3333333333333333
$data[14,6] = 0.02581033333333333
$data[14,7] = 0.077431
$data[14,8] = 0.002360902975946389
$data[14,9] = 0.002360902975946389
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 0.9008219999999999
$data[14,13] = 2.702466
$data[14,14] = 0.1339046737357332
$data[14,15] = 0.1450544297764824
$data[14,16] = 0.023250516094
$data[14,17] = 0.209254644846
$data[14,18] = 0.0003161359427158228
$data[14,19] = 0.0003424594349335039

# Row 17: "Neutro", "Cd5l", "Cd5" -> "M2"
$data[15,0] = "Neutro"
$data[15,1] = "Cd5l"
$data[15,2] = "Cd5"
$data[15,3] = "M2"
$data[15,4] = 1
$data[15,5] = 0.3333333333333333
$data[15,6] = 0.02581033333333333
$data[15,7] = 0.077431
$data[15,8] = 0.002360902975946389
$data[15,9] = 0.002360902975946389
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 2.363391666666667
$data[15,13] = 7.090175
$data[15,14] = 0.3513115688057693
$data[15,15] = 0.3805640077027689
$data[15,16] = 0.0609999267138889
$data[15,17] = 0.548999340425
$data[15,18] = 0.0008294125282779354
$data[15,19] = 0.0008984746983235519

# Row 18: "Neutro", "Cd5l", "Cd5" -> "Neutro"
$data[16,0] = "Neutro"
$data[16,1] = "Cd5l"
$data[16,2] = "Cd5"
$data[16,3] = "Neutro"
$data[16,4] = 1
$data[16,5] = 0.3333333333333333
$data[16,6] = 0.02581033333333333
$data[16,7] = 0.077431
$data[16,8] = 0.002360902975946389
$data[16,9] = 0.002360902975946389
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 0.414697
$data[16,13] = 1.244091
$data[16,14] = 0.06164355053960423
$data[16,15] = 0.06677638519598536
$data[16,16] = 0.01070346780233333
$data[16,17] = 0.09633121022100001
$data[16,18] = 0.0001455344419168533
$data[16,19] = 0.0001576525665321443

# Row 19: "Neutro", "Cd5l", "Cd5" -> "sCs"
$data[17,0] = "Neutro"
$data[17,1] = "Cd5l"
$data[17,2] = "Cd5"
$data[17,3] = "sCs"
$data[17,4] = 1
$data[17,5] = 0.3333333333333333
$data[17,6] = 0.02581033333333333
$data[17,7] = 0.077431
$data[17,8] = 0.002360902975946389
$data[17,9] = 0.002360902975946389
$data[17,10] = 2
$data[17,11] = 1
$data[17,12] = 1.036213
$data[17,13] = 2.072426
$data[17,14] = 0.1540301676532382
$data[17,15] = 0.1112371336712308
$data[17,16] = 0.02674500293433334
$data[17,17] = 0.160470017606
$data[17,18] = 0.0003636502811980513
$data[17,19] = 0.0002626200799201551

# Write the whole block back in one shot so the used range (and dimension) is A1:T19
$ws.Range("A2:T19").Value = $data
